# Add a new "2022" data column (N) to the worksheet, mirroring the style
# of the existing 2021 column (M), and populate it with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: N4 = 2022
$ws.Range("N4").Value = 2022

# Data rows 5-13: new 2022 values
$ws.Range("N5").Value = 4.3
$ws.Range("N6").Value = 5.1
$ws.Range("N7").Value = 3.1
$ws.Range("N8").Value = 2.9
$ws.Range("N9").Value = 3.4
$ws.Range("N10").Value = 2.3
$ws.Range("N11").Value = 92.8
$ws.Range("N12").Value = 91.6
$ws.Range("N13").Value = 94.6

# Copy styles from column M (2021) into the new column N (2022) so
# formatting (number format, borders, fonts) matches.
$ws.Range("M4:M13").Copy()
$ws.Range("N4:N13").PasteSpecial(-4122)  # xlPasteFormats

# Move the active selection to N15, matching the post-edit workbook state.
$ws.Range("N15").Select()
